$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All target cells in this sheet are stored as text (inlineStr) in the
# original workbook, including values that look numeric (e.g. prices like
# "315.46" using "." as a thousands separator, or percentages like
# "  +2.34%  "). Force a text number format before assigning so Excel
# does not auto-convert these strings into numbers/dates/percentages.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.489.70"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.828.64"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.10%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.46"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5073"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -4.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3912"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07702"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +3.82%  "
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.93"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.60%  "
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "Polygon"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.115"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.07"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.281"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.76%  "
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "BinanceUSD"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.001"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.04%  "
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.566"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.828.00"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.60"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +6.31%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06646"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.71"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.94%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.161"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +3.71%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.519.54"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.255"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +7.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.64"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.036.47"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.410"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.27"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.134"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.47%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1088"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.672"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.36%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07071"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.2224"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.891"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.153"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6252"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.189"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.51"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5909"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +4.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.716"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.91"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.981"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.85%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.195"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.31%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06924"
